# Fruta / hortaliza, semanal
# Insert two new weekly price rows (173-174) for "Macroferia Regional de
# Talca - Kiwi", pushing the existing rows 173-185 down to 175-187.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right above the current row 173 (each Insert()
# shifts row 173 and everything below it down by one).
$ws.Rows.Item(173).Insert()
$ws.Rows.Item(173).Insert()

# Populate the two new rows with the same constant descriptive columns
# used throughout this market/product block (A, B, C, E, F, G, H, I, J, K, R),
# plus the week-specific values (D, L, M, N, O, P, Q, S, T).

$newRows = @(
  @{ Row = 173; D = 44461; L = "Especial"; M = 200; N = 12000; O = 12000; P = 12000; Q = "$/bandeja 18 kilos"; S = 667; T = 18 },
  @{ Row = 174; D = 44461; L = "Primera";  M = 150; N = 10000; O = 10000; P = 10000; Q = "$/bandeja 18 kilos"; S = 556; T = 18 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 5
    $ws.Cells.Item($row, 2).Value = "Macroferia Regional de Talca"
    $ws.Cells.Item($row, 3).Value = "Maule"
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = 7
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100101
    $ws.Cells.Item($row, 8).Value = "Berries"
    $ws.Cells.Item($row, 9).Value = 100101007
    $ws.Cells.Item($row, 10).Value = "Kiwi"
    $ws.Cells.Item($row, 11).Value = "Hayward"
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = "Provincia de Curicó"
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
